$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Remove the duplicate "Contact" row (row 11) entirely - this row was a
# leftover duplicate of row 10 ("Contact" / "No display for ContactDetail").
# Deleting it shifts rows 12-22 up by one, matching the new A1:B21 dimension.
$ws.Rows.Item(11).Delete()

# Version bump
$ws.Range("B3").Value = "6.0.0"

# Date bump
$ws.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value (was previously blank)
$ws.Range("B9").Value = "Alvearie Team"

# The old "Contact" row (row 10) becomes "Jurisdiction"
$ws.Range("A10").Value = "Jurisdiction"
$ws.Range("B10").Value = "United States of America"

# Case Sensitive value. Assigning the literal text "true" through .Value
# gets auto-coerced to an Excel boolean TRUE, so instead write a formula
# that evaluates to the text "true" and then convert it in place to a
# plain value via copy / paste-special-values (keeps it a real text
# string and keeps the existing cell style untouched).
$ws.Range("B14").Formula = '="true"'
$ws.Range("B14").Copy()
$ws.Range("B14").PasteSpecial(-4163)
